# tests/data/valid.mirri.xlsx — "update test files. some validation changes"
#
# 1. Rename two worksheet tabs:
#      "Sexual states" -> "Sexual state"
#      "Ontobiotype"   -> "Ontobiotope"
# 2. Move the active/selected tab from "Strains" (index 1) to the
#    renamed "Sexual state" sheet (index 2), which flips each sheet's
#    tabSelected flag and the workbook's bookViews/activeTab.
# 3. Re-affirm the (already-collapsed) single-cell selection on every
#    sheet so the saved sqref/activeCell reflect just the active cell
#    rather than the old multi-area selection.

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheets -------------------------------------------------
$wb.Worksheets.Item("Sexual states").Name = "Sexual state"
$wb.Worksheets.Item("Ontobiotype").Name = "Ontobiotope"

# --- 3. Restore single-cell selections on each sheet -------------------
$null = $wb.Worksheets.Item("Geographic origin").Range("E79").Select()
$null = $wb.Worksheets.Item("Strains").Range("AE21").Select()
$null = $wb.Worksheets.Item("Sexual state").Range("A1").Select()
$null = $wb.Worksheets.Item("Genomic information").Range("A6").Select()
$null = $wb.Worksheets.Item("Markers").Range("A1").Select()
$null = $wb.Worksheets.Item("Literature").Range("A1").Select()
$null = $wb.Worksheets.Item("Growth media").Range("A12").Select()
$null = $wb.Worksheets.Item("Ontobiotope").Range("A1").Select()

# --- 2. Make "Sexual state" the active tab (was "Strains") -------------
$wb.Worksheets.Item("Sexual state").Activate()
